$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.149.21"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.583.97"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.47%  "
$ws.Range("D9").Value = "2.591.31"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.19%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "3.043.72"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").Value = "60.256.51"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("D18").Value = "2.595.75"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "2.703.49"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "0.0₃0844"
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +17.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "296.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.615"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0560"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("E51").Value = "  -0.05%  "
